$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.772.10'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.536.63'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '205.28'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '21.26'
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").Value = '1.754.40'
$ws.Range("E12").Value = '  -1.84%  '
$ws.Range("D13").Value = '1.535.97'
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("D16").Value = '26.763.70'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '60.95'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '212.79'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").Value = '0.0₃0681'
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '4.01'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("E24").Value = '  -3.48%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '151.27'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  -2.38%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '14.78'
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").Value = '1.362.33'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.961'
$ws.Range("E36").Value = '  +4.27%  '
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.802'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '5.73'
$ws.Range("E41").Value = '  +7.52%  '
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '62.73'
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("E45").Value = '  -3.53%  '
$ws.Range("D46").Value = '1.669.57'
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '85.13'
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.0506'
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("D49").Value = '0.0₇0974'
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("E51").Value = '  -0.16%  '
